$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 270-301 (data through 2021-06-28)
$data = @(
    ,@(44344, 1, 6, 35.01604902246864)
    ,@(44345, 0, 6, 35.01604902246864)
    ,@(44346, 1, 6, 35.01604902246864)
    ,@(44347, 3, 9, 52.52407353370295)
    ,@(44348, 0, 9, 52.52407353370295)
    ,@(44349, 2, 7, 40.85205719288007)
    ,@(44350, 0, 7, 40.85205719288007)
    ,@(44351, 0, 6, 35.01604902246864)
    ,@(44352, 0, 6, 35.01604902246864)
    ,@(44353, 1, 6, 35.01604902246864)
    ,@(44354, 0, 3, 17.50802451123432)
    ,@(44355, 0, 3, 17.50802451123432)
    ,@(44356, 1, 2, 11.67201634082288)
    ,@(44357, 1, 3, 17.50802451123432)
    ,@(44358, 0, 3, 17.50802451123432)
    ,@(44359, 0, 3, 17.50802451123432)
    ,@(44360, 1, 3, 17.50802451123432)
    ,@(44361, 0, 3, 17.50802451123432)
    ,@(44362, 2, 5, 29.18004085205719)
    ,@(44363, 1, 5, 29.18004085205719)
    ,@(44364, 5, 9, 52.52407353370295)
    ,@(44365, 0, 9, 52.52407353370295)
    ,@(44366, 0, 9, 52.52407353370295)
    ,@(44367, 0, 8, 46.68806536329151)
    ,@(44368, 0, 8, 46.68806536329151)
    ,@(44369, 0, 6, 35.01604902246864)
    ,@(44370, 0, 5, 29.18004085205719)
    ,@(44371, 0, 0, 0)
    ,@(44372, 0, 0, 0)
    ,@(44373, 0, 0, 0)
    ,@(44374, 0, 0, 0)
    ,@(44375, 0, 0, 0)
)

$startRow = 270
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Apply the same date style (format YYYY-MM-DD HH:MM:SS, s="2") used by the
# rest of column A to the newly added date cells, without introducing a new style.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
